# Tensile test mapping workbook: introduce dcterms-style IRI mapping,
# clear out the per-row "Annotation" placeholder values (keep the one
# genuine ontology IRI on the Werkstoff/Material row), widen the
# Class-type / Annotation columns to fit the new long IRIs, and recolor
# the header font from a theme color to an explicit black.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Column B ("Class type"): replace the short local names with the
#    fully qualified https://w3id.org/steel/ProcessOntology/... IRIs.
#    (B2 "TestingFacility" is left as-is -- it was not touched upstream.)
# ---------------------------------------------------------------------
$classTypeMap = @{
  3  = "https://w3id.org/steel/ProcessOntology/ProjectNumber"
  4  = "https://w3id.org/steel/ProcessOntology/ProjectName"
  5  = "https://w3id.org/steel/ProcessOntology/TimeStamp"
  6  = "https://w3id.org/steel/ProcessOntology/MachineData"
  7  = "https://w3id.org/steel/ProcessOntology/ForceMeasuringDevice"
  8  = "https://w3id.org/steel/ProcessOntology/DisplacementTransducer"
  9  = "https://w3id.org/steel/ProcessOntology/TestStandard"
  10 = "https://w3id.org/steel/ProcessOntology/Material"
  11 = "https://w3id.org/steel/ProcessOntology/SpecimenType"
  12 = "https://w3id.org/steel/ProcessOntology/Tester"
  13 = "https://w3id.org/steel/ProcessOntology/SampleIdentifier-2"
  14 = "https://w3id.org/steel/ProcessOntology/OriginalGaugeLength"
  15 = "https://w3id.org/steel/ProcessOntology/ParallelLength"
  16 = "https://w3id.org/steel/ProcessOntology/SpecimenThickness"
  17 = "https://w3id.org/steel/ProcessOntology/SpecimenWidth"
  18 = "https://w3id.org/steel/ProcessOntology/TestingRate"
  19 = "https://w3id.org/steel/ProcessOntology/Preload"
  20 = "https://w3id.org/steel/ProcessOntology/Temperature"
  21 = "https://w3id.org/steel/ProcessOntology/Remark"
  22 = "https://w3id.org/steel/ProcessOntology/TestTime"
  23 = "https://w3id.org/steel/ProcessOntology/StandardForce"
  24 = "https://w3id.org/steel/ProcessOntology/AbsoluteCrossheadTravel"
  25 = "https://w3id.org/steel/ProcessOntology/Extension"
  26 = "https://w3id.org/steel/ProcessOntology/WidthChange"
  27 = "https://w3id.org/steel/ProcessOntology/PercentageElongation"
}

foreach ($row in $classTypeMap.Keys) {
  $ws.Range("B$row").Value2 = $classTypeMap[$row]
}

# ---------------------------------------------------------------------
# 2. Column C ("Annotation"): these were all placeholder empty strings;
#    clear them out entirely (row 10, Werkstoff/Material, keeps its
#    real ontology-class IRI and is left untouched). Also clear the
#    stray quote-prefixed empty string left in A28.
#    Copy/PasteSpecial(formats) re-uses the plain "general alignment"
#    style instead of the old quote-prefixed one.
# ---------------------------------------------------------------------
$ws.Range("B2").Copy() | Out-Null
for ($row = 2; $row -le 27; $row++) {
  if ($row -ne 10) {
    $ws.Range("C$row").PasteSpecial(-4122) | Out-Null
  }
}
$ws.Range("B28").Copy() | Out-Null
$ws.Range("A28").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

for ($row = 2; $row -le 27; $row++) {
  if ($row -ne 10) {
    $ws.Range("C$row").ClearContents()
  }
}
$ws.Range("A28").ClearContents()

# ---------------------------------------------------------------------
# 3. Row 2 grows a little taller to match the new header block, and
#    columns B/C widen to fit the long IRI strings now stored there.
# ---------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 20.25
$ws.Columns.Item(2).ColumnWidth = 62.028809523809524
$ws.Columns.Item(3).ColumnWidth = 63.600238095238095

# ---------------------------------------------------------------------
# 4. Header font: was a theme-based color, now an explicit black.
#    (Applies to every cell that shares the bold/bordered header style,
#    i.e. the header row plus A2.)
# ---------------------------------------------------------------------
$ws.Range("A1:C1,A2").Font.Color = 0
